# Update stats for 2025-11 (row 24)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = 6364
$ws.Range("D24").Value = 5958495
$ws.Range("E24").Value = 936.2814267756128
$ws.Range("F24").Value = 8.489601091033073
$ws.Range("H24").Value = 26.22400420965933
